# Fix table for perc_cis function: update the confidence-interval
# percentage strings in the "Percentages of drivers testing positive by
# drug type, sex, and year group" table (F 2007-2010 / M 2007-2010
# columns) to the corrected values.

$d = $word.ActiveDocument

function Replace-Exact($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

# Alcohol row
Replace-Exact "38.3% (35.2%, 41.4%)" "38.7% (36.5%, 40.9%)"
Replace-Exact "39.8% (36.4%, 43.2%)" "39.1% (36.7%, 41.5%)"

# Cannabinoid row
Replace-Exact "4.6% (3.2%, 6%)" "4.1% (3.2%, 5%)"
Replace-Exact "12.2% (9.9%, 14.5%)" "11.8% (10.2%, 13.4%)"

# Depressant row
Replace-Exact "1.9% (1%, 2.8%)" "2.1% (1.4%, 2.8%)"
Replace-Exact "5.4% (3.8%, 7%)" "4.9% (3.8%, 6%)"

# Narcotic row
Replace-Exact "2% (1.1%, 2.9%)" "1.8% (1.2%, 2.4%)"
Replace-Exact "4.9% (3.4%, 6.4%)" "5.1% (4%, 6.2%)"

# Other row
Replace-Exact "4.3% (3%, 5.6%)" "3.7% (2.8%, 4.6%)"
Replace-Exact "5.1% (3.6%, 6.6%)" "5.2% (4.1%, 6.3%)"

# Stimulant row
Replace-Exact "8.1% (6.3%, 9.9%)" "8.1% (6.8%, 9.4%)"
Replace-Exact "8.9% (6.9%, 10.9%)" "8.6% (7.2%, 10%)"
